# "housekeeping and updating AoN and fine lsp results"
#
# 1) opus_base Simple aWCE: H15 -> "Loaded"
# 2) opus_base AoN aWCE: H3 -> "Loaded"
# 3) opus_base Fine aWCE: H14 -> "Loaded"
# 4) opus_base LSP AoN aWCE: fill in several previously-blank
#    Compute/Runtime/BLEU rows (C/D/F, with E recalculating from D)
# 5) opus_base LSP Fine aWCE: fill in rows 5-7 (C/D/F, with E recalculating from D)

$wb = $excel.ActiveWorkbook

# --- 1) Simple aWCE ---------------------------------------------------
$wsSimple = $wb.Worksheets.Item("opus_base Simple aWCE")
$wsSimple.Range("H15").Value = "Loaded"

# --- 2) AoN aWCE --------------------------------------------------------
$wsAoN = $wb.Worksheets.Item("opus_base AoN aWCE")
$wsAoN.Range("H3").Value = "Loaded"

# --- 3) Fine aWCE ---------------------------------------------------------
$wsFine = $wb.Worksheets.Item("opus_base Fine aWCE")
$wsFine.Range("H14").Value = "Loaded"

# --- 4) LSP AoN aWCE ------------------------------------------------------
$wsLspAoN = $wb.Worksheets.Item("opus_base LSP AoN aWCE ")

$wsLspAoN.Range("C2").Value = 42.400700000000001
$wsLspAoN.Range("D2").Value = 15672.7017
$wsLspAoN.Range("F2").Value = 69.000799999999998

$wsLspAoN.Range("C3").Value = 42.695300000000003
$wsLspAoN.Range("D3").Value = 23496.607100000001
$wsLspAoN.Range("F3").Value = 137.96010000000001

$wsLspAoN.Range("C5").Value = 42.723500000000001
$wsLspAoN.Range("D5").Value = 16049.419900000001
$wsLspAoN.Range("F5").Value = 82.794399999999996

$wsLspAoN.Range("C8").Value = 42.625599999999999
$wsLspAoN.Range("D8").Value = 17697.437999999998
$wsLspAoN.Range("F8").Value = 117.2388

$wsLspAoN.Range("C9").Value = 42.5426
$wsLspAoN.Range("D9").Value = 17531.4149
$wsLspAoN.Range("F9").Value = 117.2388

$wsLspAoN.Range("C10").Value = 42.702599999999997
$wsLspAoN.Range("D10").Value = 17335.9588
$wsLspAoN.Range("F10").Value = 82.794399999999996

$wsLspAoN.Range("C11").Value = 42.674900000000001
$wsLspAoN.Range("D11").Value = 15932.497300000001
$wsLspAoN.Range("F11").Value = 82.794399999999996

$wsLspAoN.Range("C15").Value = 42.1188
$wsLspAoN.Range("D15").Value = 13728.5939
$wsLspAoN.Range("F15").Value = 82.794399999999996

$wsLspAoN.Range("C16").Value = 42.131500000000003
$wsLspAoN.Range("D16").Value = 13429.059300000001
$wsLspAoN.Range("F16").Value = 55.2

$wsLspAoN.Range("C17").Value = 42.167900000000003
$wsLspAoN.Range("D17").Value = 15779.766900000001
$wsLspAoN.Range("F17").Value = 69.000799999999998

# --- 5) LSP Fine aWCE -------------------------------------------------------
$wsLspFine = $wb.Worksheets.Item("opus_base LSP Fine aWCE ")

$wsLspFine.Range("C5").Value = 42.472499999999997
$wsLspFine.Range("D5").Value = 19031.651600000001
$wsLspFine.Range("F5").Value = 131.08879999999999

$wsLspFine.Range("C6").Value = 42.544699999999999
$wsLspFine.Range("D6").Value = 15466.3519
$wsLspFine.Range("F6").Value = 103.4885

$wsLspFine.Range("C7").Value = 42.8215
$wsLspFine.Range("D7").Value = 25892.841199999999
$wsLspFine.Range("F7").Value = 137.96010000000001

# --- selection bookkeeping (matches each sheet's saved cursor position) ---
$wsSimple.Range("I20").Select()
$wsAoN.Range("H11").Select()
$wsFine.Range("I21").Select()
$wsLspAoN.Range("D30").Select()
# LSP Fine aWCE is the tab that was active/selected when the file was saved,
# so leave it selected last.
$wsLspFine.Range("C28").Select()
